$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '27.372.02'
Set-TextValue 'E2' '  -3.86%  '

# Row 3
Set-TextValue 'D3' '1.860.83'
Set-TextValue 'E3' '  -4.63%  '

# Row 4
Set-TextValue 'D4' '1.001'
Set-TextValue 'E4' '  -1.23%  '

# Row 5
Set-TextValue 'D5' '323.41'
Set-TextValue 'E5' '  +0.40%  '

# Row 6
Set-TextValue 'D6' '1.001'
Set-TextValue 'E6' '  -1.11%  '

# Row 7
Set-TextValue 'D7' '0.4539'
Set-TextValue 'E7' '  -5.27%  '

# Row 8
Set-TextValue 'D8' '0.3873'
Set-TextValue 'E8' '  -5.11%  '

# Row 9
Set-TextValue 'D9' '48.40'
Set-TextValue 'E9' '  -10.34%  '

# Row 10
Set-TextValue 'D10' '0.07915'

# Row 11
Set-TextValue 'E11' '  -3.27%  '

# Row 12
Set-TextValue 'D12' '21.47'
Set-TextValue 'E12' '  -4.26%  '

# Row 13
Set-TextValue 'D13' '1.854.44'
Set-TextValue 'E13' '  -5.52%  '

# Row 14
Set-TextValue 'D14' '5.917'
Set-TextValue 'E14' '  -3.97%  '

# Row 15
Set-TextValue 'D15' '7.137'
Set-TextValue 'E15' '  -5.72%  '

# Row 16
Set-TextValue 'D16' '1.001'
Set-TextValue 'E16' '  -1.35%  '

# Row 17
Set-TextValue 'D17' '0.00001035'
Set-TextValue 'E17' '  -3.49%  '

# Row 18
Set-TextValue 'D18' '85.85'
Set-TextValue 'E18' '  -5.18%  '

# Row 19
Set-TextValue 'D19' '0.06525'
Set-TextValue 'E19' '  -1.65%  '

# Row 20
Set-TextValue 'D20' '17.16'
Set-TextValue 'E20' '  -6.99%  '

# Row 21
Set-TextValue 'D21' '1.001'
Set-TextValue 'E21' '  -1.05%  '

# Row 22
Set-TextValue 'D22' '5.537'
Set-TextValue 'E22' '  -5.22%  '

# Row 23
Set-TextValue 'D23' '27.372.30'
Set-TextValue 'E23' '  -3.92%  '

# Row 24
Set-TextValue 'D24' '10.88'
Set-TextValue 'E24' '  -5.00%  '

# Row 25
Set-TextValue 'E25' '  -1.15%  '

# Row 26
Set-TextValue 'D26' '2.066.90'
Set-TextValue 'E26' '  -5.80%  '

# Row 27
Set-TextValue 'D27' '153.04'
Set-TextValue 'E27' '  -2.19%  '

# Row 28
Set-TextValue 'D28' '19.79'
Set-TextValue 'E28' '  -2.56%  '

# Row 29
Set-TextValue 'D29' '2.068'
Set-TextValue 'E29' '  -4.85%  '

# Row 30
Set-TextValue 'D30' '5.486'
Set-TextValue 'E30' '  -5.69%  '

# Row 31
Set-TextValue 'D31' '120.81'
Set-TextValue 'E31' '  -2.82%  '

# Row 32
Set-TextValue 'D32' '1.488'
Set-TextValue 'E32' '  +3.14%  '

# Row 33
Set-TextValue 'D33' '0.09323'
Set-TextValue 'E33' '  -3.57%  '

# Row 34
Set-TextValue 'D34' '0.9352'
Set-TextValue 'E34' '  -5.03%  '

# Row 35
Set-TextValue 'D35' '3.614'
Set-TextValue 'E35' '  -2.18%  '

# Row 36
Set-TextValue 'D36' '5.273'
Set-TextValue 'E36' '  -6.28%  '

# Row 37
Set-TextValue 'D37' '0.02238'
Set-TextValue 'E37' '  -4.11%  '

# Row 40
Set-TextValue 'D40' '8.266'
Set-TextValue 'E40' '  -9.26%  '

# Row 41
Set-TextValue 'D41' '1.000'
Set-TextValue 'E41' '  -1.13%  '

# Row 42
Set-TextValue 'D42' '0.5923'
Set-TextValue 'E42' '  -4.83%  '

# Row 43
Set-TextValue 'D43' '0.1889'
Set-TextValue 'E43' '  -1.61%  '

# Row 44
Set-TextValue 'D44' '10.13'
Set-TextValue 'E44' '  -9.58%  '

# Row 45
Set-TextValue 'D45' '1.262'
Set-TextValue 'E45' '  -5.67%  '

# Row 46
Set-TextValue 'D46' '0.5639'
Set-TextValue 'E46' '  -5.22%  '

# Row 47
Set-TextValue 'D47' '11.94'
Set-TextValue 'E47' '  -7.47%  '

# Row 50
Set-TextValue 'D50' '0.06783'
Set-TextValue 'E50' '  -0.51%  '

# Row 51
Set-TextValue 'D51' '108.71'
Set-TextValue 'E51' '  -1.98%  '

# Row 38 (swap)
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D38' '1.223'
Set-TextValue 'E38' '  -2.16%  '

# Row 39 (swap)
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D39' '0.05993'
Set-TextValue 'E39' '  -3.23%  '

# Row 48 (swap)
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D48' '3.372'
Set-TextValue 'E48' '  -1.13%  '

# Row 49 (swap)
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D49' '1.926'
Set-TextValue 'E49' '  -6.48%  '
